# Applies the workbook update:
#   1. Column C ("Förändrad") date value changes from 45184 to 45186 for every data row.
#   2. Existing HYPERLINK(...) formulas (columns S,T,U,V,W,X,Y) get a second
#      "friendly name" argument equal to the row's Beteckning (column A) value,
#      e.g. HYPERLINK("url") -> HYPERLINK("url", "A 45074-2020")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (header is row 1, data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

# 1. Bulk-update column C for all data rows in one shot (value only, keeps the
#    existing date number format / style untouched).
$ws.Range("C2:C" + $lastRow).Value = 45186

# 2. Walk every data row and patch any HYPERLINK formulas that still only have
#    a single (url) argument, adding the row's column-A text as display name.
$cols = @("S", "T", "U", "V", "W", "X", "Y")
for ($r = 2; $r -le $lastRow; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($aVal)) { continue }

    foreach ($col in $cols) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) { continue }
        if ($f -notmatch '^=HYPERLINK\(') { continue }
        if ($f.TrimEnd() -match ',\s*"[^"]*"\s*\)$') { continue }

        $newF = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
        $cell.Formula = $newF
    }
}

Write-Host "Updated column C and hyperlink formulas through row $lastRow"
